# Update price/CO2 scenario table: fill in the forecast values for
# 2021-2051 (columns P:AT) on rows 2 and 3, replacing the placeholder
# zeros copied from 2020, and apply a plain integer number format to
# those newly-populated cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @(
    25, 30, 35, 45, 55, 57, 60, 62, 64, 65, 67, 71, 75, 80, 84,
    83.808348127145933, 86.158059642158435, 88.433730600574123,
    90.641623163612849, 92.787171856914142, 94.875128470163503,
    96.909675952759088, 98.894519037239718, 100.83295714376011,
    102.72794361897078, 104.58213431134067, 106.39792773605862,
    108.17749854124806, 109.92282559050837, 111.63571568246282, 115
)

# Columns P (16) through AT (46) inclusive -> 31 columns, for both data rows.
$startCol = 16
$endCol = 46

for ($col = $startCol; $col -le $endCol; $col++) {
    $value = $values[$col - $startCol]
    $ws.Cells.Item(2, $col).Value = $value
    $ws.Cells.Item(3, $col).Value = $value
}

# The refreshed forecast years get a plain integer display format.
$rng = $ws.Range($ws.Cells.Item(2, $startCol), $ws.Cells.Item(3, $endCol))
$rng.NumberFormat = "0"

# Leave the sheet scrolled over to the newly-edited columns with the
# last column of both rows selected, matching where the edit was made.
[void]$ws.Activate()
[void]$ws.Range("AB1").Select()
[void]$ws.Range("AT2:AT3").Select()
